$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '69.851.81'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -1.95%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.760.00'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +2.33%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '620.79'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +2.86%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '182.01'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.04%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.757.54'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +2.39%  '
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.04%  '
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -0.37%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.168'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +3.30%  '
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -4.97%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.492'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -1.42%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '41.56'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +2.27%  '
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +1.95%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.378.00'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +2.14%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.753.90'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +2.23%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '69.951.82'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -1.74%  '
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +0.23%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.61'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +1.34%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '16.79'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -1.72%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '508.84'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -2.21%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.57'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +3.60%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.729'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -2.13%  '
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +1.65%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '87.35'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -1.00%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '13.18'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -2.60%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.14'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +1.27%  '
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +20.19%  '
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +0.11%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.54'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -0.61%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.92'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +4.42%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.96'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -2.14%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '31.13'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -2.09%  '
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -0.84%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.999'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -0.02%  '
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +4.41%  '
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +0.89%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.338'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -2.42%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.133'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +2.58%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.12'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -3.05%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '50.28'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -2.24%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '45.62'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +0.77%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '428.40'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +3.13%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.76'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -0.71%  '
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +1.99%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.010.86'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -3.91%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0366'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -0.77%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '27.54'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -3.20%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '137.61'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -1.56%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.52'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +1.90%  '
